$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 44446
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("P2").Value = 467
$ws.Range("D3").Value = 44453
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("P3").Value = 400
$ws.Range("D5").Value = 44421
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15400
$ws.Range("P5").Value = 513
$ws.Range("D6").Value = 44432
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("O6").Value = 'Provincia del Elquí'
$ws.Range("P6").Value = 467
$ws.Range("D7").Value = 44376
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 18000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 18000
$ws.Range("P7").Value = 600
$ws.Range("D8").Value = 44467
$ws.Range("J8").Value = 35
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("P8").Value = 400
$ws.Range("D9").Value = 44841
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("O9").Value = 'Provincia de Limarí'
$ws.Range("P9").Value = 400
$ws.Range("D10").Value = 44841
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("N10").Value = '$/caja 40 unidades'
$ws.Range("P10").Value = 250
$ws.Range("Q10").Value = 40
$ws.Range("D11").Value = 44474
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 333
$ws.Range("D12").Value = 44418
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = '$/caja 30 unidades'
$ws.Range("P12").Value = 500
$ws.Range("Q12").Value = 30
$ws.Range("D15").Value = 44425
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 14000
$ws.Range("P15").Value = 467
$ws.Range("D16").Value = 44460
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 13000
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = 13000
$ws.Range("P16").Value = 433
